# Automated data refresh for the MeteoCat daily summary sheet.
# Updates the extraction timestamp (col E) and the refreshed observation
# values for each station row, matching commit "Update automàtic: dades i
# banners [2026-02-28 20:50]".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-28 20:48:33'
$ws.Range('I2').Value = '0.3 mm'
$ws.Range('K2').Value = '11.7 MJ/m2'
$ws.Range('E3').Value = '2026-02-28 20:48:36'
$ws.Range('N3').Value = '-3.1 °C 20:28 TU'
$ws.Range('O3').Value = '-1.3 °C'
$ws.Range('E4').Value = '2026-02-28 20:48:38'
$ws.Range('H4').Value = '''83%'
$ws.Range('O4').Value = '11.0 °C'
$ws.Range('E5').Value = '2026-02-28 20:48:41'
$ws.Range('N5').Value = '-3.0 °C 20:28 TU'
$ws.Range('E6').Value = '2026-02-28 20:48:44'
$ws.Range('N6').Value = '9.3 °C 20:29 TU'
$ws.Range('E7').Value = '2026-02-28 20:48:46'
$ws.Range('E8').Value = '2026-02-28 20:48:49'
$ws.Range('H8').Value = '''92%'
$ws.Range('E9').Value = '2026-02-28 20:48:51'
$ws.Range('E10').Value = '2026-02-28 20:48:54'
$ws.Range('H10').Value = '''86%'
$ws.Range('O10').Value = '10.8 °C'
$ws.Range('E11').Value = '2026-02-28 20:48:57'
$ws.Range('O11').Value = '7.3 °C'
$ws.Range('E12').Value = '2026-02-28 20:48:59'
$ws.Range('E13').Value = '2026-02-28 20:49:02'
$ws.Range('J13').Value = '1024.2 hPa'
$ws.Range('E14').Value = '2026-02-28 20:49:05'
$ws.Range('E15').Value = '2026-02-28 20:49:07'
$ws.Range('O15').Value = '11.0 °C'
$ws.Range('E16').Value = '2026-02-28 20:49:10'
$ws.Range('H16').Value = '''67%'
$ws.Range('N16').Value = '-3.1 °C 20:29 TU'
$ws.Range('O16').Value = '-1.2 °C'
$ws.Range('E17').Value = '2026-02-28 20:49:12'
$ws.Range('N17').Value = '1.4 °C 20:24 TU'
$ws.Range('O17').Value = '2.9 °C'
$ws.Range('E18').Value = '2026-02-28 20:49:15'
$ws.Range('H18').Value = '''83%'
$ws.Range('O18').Value = '11.6 °C'
$ws.Range('E19').Value = '2026-02-28 20:49:18'
$ws.Range('H19').Value = '''78%'
$ws.Range('E20').Value = '2026-02-28 20:49:20'
$ws.Range('H20').Value = '''64%'
$ws.Range('N20').Value = '-2.2 °C 20:19 TU'
$ws.Range('E21').Value = '2026-02-28 20:49:23'
$ws.Range('O21').Value = '7.7 °C'
$ws.Range('E22').Value = '2026-02-28 20:49:25'
$ws.Range('H22').Value = '''69%'
$ws.Range('E23').Value = '2026-02-28 20:49:28'
$ws.Range('H23').Value = '''71%'
$ws.Range('I23').Value = '1.5 mm'
$ws.Range('N23').Value = '-3.0 °C 20:22 TU'
$ws.Range('O23').Value = '-0.4 °C'
$ws.Range('E24').Value = '2026-02-28 20:49:31'
$ws.Range('J24').Value = '1025.2 hPa'
$ws.Range('E25').Value = '2026-02-28 20:49:33'
$ws.Range('H25').Value = '''63%'
$ws.Range('N25').Value = '-1.1 °C 20:26 TU'
$ws.Range('O25').Value = '1.2 °C'
$ws.Range('E26').Value = '2026-02-28 20:49:36'
$ws.Range('H26').Value = '''81%'
$ws.Range('J26').Value = '1024.4 hPa'
$ws.Range('O26').Value = '4.8 °C'
$ws.Range('E27').Value = '2026-02-28 20:49:39'
$ws.Range('H27').Value = '''56%'
$ws.Range('N27').Value = '-0.6 °C 20:28 TU'
$ws.Range('O27').Value = '1.8 °C'
$ws.Range('E28').Value = '2026-02-28 20:49:41'
$ws.Range('E29').Value = '2026-02-28 20:49:44'
$ws.Range('O29').Value = '11.7 °C'
$ws.Range('E30').Value = '2026-02-28 20:49:47'
$ws.Range('J30').Value = '1024.8 hPa'
$ws.Range('E31').Value = '2026-02-28 20:49:49'
$ws.Range('L31').Value = '66.6 km/h - 347º 20:03 TU'
$ws.Range('E32').Value = '2026-02-28 20:49:52'
$ws.Range('E33').Value = '2026-02-28 20:49:55'
$ws.Range('J33').Value = '1023.1 hPa'
$ws.Range('O33').Value = '7.1 °C'
$ws.Range('E34').Value = '2026-02-28 20:49:57'
$ws.Range('I34').Value = '0.7 mm'
$ws.Range('E35').Value = '2026-02-28 20:50:00'
$ws.Range('H35').Value = '''84%'
$ws.Range('J35').Value = '1024.9 hPa'
$ws.Range('E36').Value = '2026-02-28 20:50:02'
$ws.Range('H36').Value = '''79%'
$ws.Range('O36').Value = '12.7 °C'
$ws.Range('E37').Value = '2026-02-28 20:50:05'
$ws.Range('J37').Value = '1025.8 hPa'
$ws.Range('O37').Value = '7.2 °C'
$ws.Range('E38').Value = '2026-02-28 20:50:08'
$ws.Range('E39').Value = '2026-02-28 20:50:10'
$ws.Range('H39').Value = '''64%'
$ws.Range('N39').Value = '-2.0 °C 20:02 TU'
$ws.Range('O39').Value = '-0.6 °C'
$ws.Range('E40').Value = '2026-02-28 20:50:13'
$ws.Range('O40').Value = '7.7 °C'
$ws.Range('E41').Value = '2026-02-28 20:50:15'
$ws.Range('J41').Value = '1024.6 hPa'
$ws.Range('E42').Value = '2026-02-28 20:50:18'
$ws.Range('O42').Value = '11.1 °C'
$ws.Range('E43').Value = '2026-02-28 20:50:20'
$ws.Range('E44').Value = '2026-02-28 20:50:23'
$ws.Range('N44').Value = '-3.0 °C 20:05 TU'
$ws.Range('E45').Value = '2026-02-28 20:50:25'
$ws.Range('H45').Value = '''85%'
$ws.Range('J45').Value = '1025.5 hPa'
$ws.Range('N45').Value = '4.2 °C 20:28 TU'
$ws.Range('O45').Value = '6.2 °C'
$ws.Range('E46').Value = '2026-02-28 20:50:28'
$ws.Range('J46').Value = '1025.1 hPa'
